# Commit message: "add . to each line"
# Append a full stop ("。") to the end of each line of the test-step
# instructions on the "summary" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("summary")

$ws.Range("B4").Value = "1.在服务器上上传最新版本的测试包或者提高版本。"
$ws.Range("B5").Value = "2.安装低版本进行升级测试。"
$ws.Range("B6").Value = "3.主要测试方向：" + [char]10 + "A.查看是否弹出升级提示。" + [char]10 + "B.如果启动画面有修改升级后需要对最新画面全部展示。" + [char]10 + "C.升级后必须依然是已经登录状态。" + [char]10 + "D.升级后新增/修改的功能是否能够正确实现。"

$ws.Range("B6").Select()
